# Update the "Förändrad" (column C) date for rows 2-42 from 2025-03-05
# (serial 45721) to 2025-03-06 (serial 45722), as recorded by the
# automatic update-of-files job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45721) {
        $cell.Value2 = 45722
    }
}
